# Add 2022-Q3 data
#
# The workbook has a "总计" (Total) summary sheet followed by one sheet per
# quarter (newest first). This edit inserts a brand-new "2022-Q3" quarter
# sheet (with its own fund-holding detail table) right after "总计" and
# before the existing "2022-Q2" sheet, and adds the corresponding summary
# row to the "总计" sheet. All the other quarter sheets keep their data
# unchanged; they just shift one tab to the right to make room.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q3" sheet by duplicating "2022-Q2" (so it keeps
#    the exact same headers / column layout / cell styles), then put it
#    in the desired position, rename it and overwrite its data values.
# ---------------------------------------------------------------------
$sourceQ2 = $wb.Worksheets.Item("2022-Q2")
$sourceQ2.Copy($sourceQ2, $null)          # new copy placed right before "2022-Q2"

$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# Fund row 1: 290012 泰信行业精选灵活配置混合A
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.75"
$newSheet.Range("E2").Value = "91.96"
$newSheet.Range("F2").Value = "5.39"
$newSheet.Range("G2").Value = "0.0404"
$newSheet.Range("H2").Value = 7

# Fund row 2: 002583 泰信行业精选灵活配置混合C
$newSheet.Range("D3:G3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.04"
$newSheet.Range("E3").Value = "91.96"
$newSheet.Range("F3").Value = "5.39"
$newSheet.Range("G3").Value = "0.0022"
$newSheet.Range("H3").Value = 7

# The NumberFormat="@" calls above stamp a (harmless but extraneous) text
# number-format style onto D2:G3. Scrub it back to the plain/default style
# those cells had originally by pasting the formatting of an already
# "clean" text cell (B2, which is plain inline/shared text with no special
# style) on top of them - this keeps the values, only resets the style.
$newSheet.Range("B2").Copy()
$newSheet.Range("D2:G3").PasteSpecial(-4122)
$newSheet.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Update the "总计" (Total) sheet: shift the existing 7 data rows
#    (rows 2-8) down by one row to rows 3-9, preserving their styles,
#    then write the new 2022-Q3 row into row 2 and fix up the running
#    index numbers in column A (0..7).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D8").Copy($total.Range("A3:D9"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
$total.Range("A9").Value = 7

$total.Application.CutCopyMode = $false
